$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case row appended at the bottom of the "Test Cases" sheet
# (TCID | JIRA ID | Description | Runmode | Results) -> row 58.
$row = 58

# --- Values -------------------------------------------------------------
$ws.Range("A$row").Value = "Profile57"
$ws.Range("B$row").Value = "OPQA-2101|OPQA-2100"
$ws.Range("C$row").Value = " `nVerify that profile modal displays the profile picture of the user.|Verify that profile modal displays text that explains introduces the profile concept."
$ws.Range("D$row").Value = "Y"

# --- Formatting: copy from the matching columns of the row above so the
# new row visually matches the rest of the table -------------------------
$ws.Range("A57").Copy() | Out-Null
$ws.Range("A$row").PasteSpecial(-4122) | Out-Null

$ws.Range("B21").Copy() | Out-Null
$ws.Range("B$row").PasteSpecial(-4122) | Out-Null

$ws.Range("C55").Copy() | Out-Null
$ws.Range("C$row").PasteSpecial(-4122) | Out-Null

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D$row").PasteSpecial(-4122) | Out-Null

$ws.Range("E57").Copy() | Out-Null
$ws.Range("E$row").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Row height matches the other wrapped / multi-part rows (e.g. row 57)
$ws.Rows.Item($row).RowHeight = 30

# --- Selection ------------------------------------------------------------
$null = $ws.Range("C38").Select()
